$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Primary Focus" bullet: drop the _GoBack bookmark that wrapped "Focus"
#    and extend the sentence with " & Project Management".
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$r = $d.Content
$r.Find.Execute(": Electromechanical Systems")
$r.Collapse(0)
$r.InsertAfter(" & Project Management")

# ---------------------------------------------------------------------------
# 2. Skills table, "Hardware/Design" column: the three bullets
#       Fabrication/Installation of Cable & Wiring
#       CAD (Solidworks/Autodesk)
#       Finite Element Analysis
#    become
#       Rapid Prototyping
#       Fabrication/Installation of Cable & Wiring
#       CAD (Solidworks/Autodesk)   <- with the _GoBack bookmark re-appearing
#                                       right after it (empty bookmark span)
# ---------------------------------------------------------------------------

# 2a. "Fabrication/Installation of Cable & Wiring" -> "Rapid Prototyping"
$r2 = $d.Content
$r2.Find.Execute("Fabrication/Installation of Cable & Wiring", $false, $false, $false, $false, $false, $true, 1, $false, "Rapid Prototyping", 2)

# 2b. "CAD (Solidworks/Autodesk)" (3 runs) -> "Fabrication/Installation of Cable & Wiring" (1 run)
$r3 = $d.Content
$r3.Find.Execute("CAD (Solidworks/Autodesk)", $false, $false, $false, $false, $false, $true, 1, $false, "Fabrication/Installation of Cable & Wiring", 2)

# 2c. "Finite Element Analysis" -> "CAD" / " " / "(Solidworks/Autodesk)" (3 runs),
#     then stamp the now-homeless _GoBack bookmark right after the new text.
$r4 = $d.Content
$r4.Find.Execute("Finite Element Analysis")
$newPara = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>CAD</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(Solidworks/Autodesk)</w:t></w:r></w:p>'
$r4.InsertXML($newPara)

$r5 = $d.Content
$r5.Find.Execute("(Solidworks/Autodesk)")
$r5.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r5)

# ---------------------------------------------------------------------------
# 3. Skills table, "Programming/Scripting" column.
#       MATLAB               -> VBA in Excel
#       LabVIEW (Former CLAD) -> MATLAB
# ---------------------------------------------------------------------------
$r6 = $d.Content
$r6.Find.Execute("MATLAB", $false, $false, $false, $false, $false, $true, 1, $false, "VBA in Excel", 2)

$r7 = $d.Content
$r7.Find.Execute("LabVIEW (Former CLAD)", $false, $false, $false, $false, $false, $true, 1, $false, "MATLAB", 2)

Write-Output "done"
